$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1").Value = " ValorPlano"
$ws.Range("A1").Select()
